# CertificationProject-HybridFramework/framework/InputData/TestData.xlsx
# Commit: "Added class mercuryFlightBooking, dataProvider."
#
# 1. Switch the browser under test from firefox to opera (openBrowser!B1).
# 2. Remove the unused, empty "Sheet1".
# 3. Add a new "inputFlightDetails" sheet (after RegisterNewUser) holding the
#    data-provider rows for the new mercuryFlightBooking test.
# 4. Move the "active tab" from RegisterNewUser back to openBrowser (first tab,
#    cell B1 selected).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# -- 1. openBrowser: firefox -> opera ---------------------------------------
$openBrowser = $wb.Worksheets.Item("openBrowser")
$openBrowser.Range("B1").Value = "opera"

# -- 2. Remove the empty placeholder sheet ----------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete() | Out-Null

# -- 3. Add the new inputFlightDetails sheet after RegisterNewUser ----------
$registerNewUser = $wb.Worksheets.Item("RegisterNewUser")
$flightSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $registerNewUser)
$flightSheet.Name = "inputFlightDetails"

$flightData = @(
    @("Journey Type", "Round Trip"),
    @("Number of Passengers", 2),
    @("Departing From", "Acapulco"),
    @("Departure Month", "September"),
    @("Departure Day", 12),
    @("Arriving At", "Zurich"),
    @("Arrival Month", "September"),
    @("Arrival Day", 13),
    @("Service Class", "Business Class"),
    @("Airline", "Blue Skies Airlines")
)

for ($i = 0; $i -lt $flightData.Count; $i++) {
    $row = $i + 1
    $flightSheet.Cells.Item($row, 1).Value = $flightData[$i][0]
    $flightSheet.Cells.Item($row, 2).Value = $flightData[$i][1]
}

$flightSheet.Columns.Item(1).ColumnWidth = 18.125

# -- 4. Restore the active tab / selection to openBrowser!B1 ----------------
$openBrowser.Activate()
$openBrowser.Range("B1").Select() | Out-Null
